$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 95, shifting existing rows 95:105 down to 96:106
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new record's data
$ws.Range("A95").Value = 5
$ws.Range("B95").Value = "Macroferia Regional de Talca"
$ws.Range("C95").Value = "Maule"
$ws.Range("D95").Value = 45194
$ws.Range("E95").Value = 7
$ws.Range("F95").Value = 300000000
$ws.Range("G95").Value = "Espárragos"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 2000
$ws.Range("K95").Value = 1200
$ws.Range("L95").Value = 1300
$ws.Range("M95").Value = 1250
$ws.Range("N95").Value = "$/kilo"
$ws.Range("O95").Value = "Provincia de Linares"
$ws.Range("P95").Value = 1250
$ws.Range("Q95").Value = 1
$ws.Range("R95").Value = "Hortaliza"
